$wb = $excel.ActiveWorkbook

# --- Teacher sheet ---
$teacher = $wb.Worksheets.Item("Teacher")

# Update existing teacher row (Marchard) with new class/unavailable-time data
$teacher.Range("B2").Value = "MATH 313, STAT 102"
$teacher.Range("C2").Value = "8,9,10"
$teacher.Range("D2").Value = "9,10"

# Add new teacher row
$teacher.Range("A3").Value = "Hurl"
$teacher.Range("B3").Value = "MATH 125"
$teacher.Range("C3").Value = "12,3,4"
$teacher.Range("D3").Value = "930,2"

$teacher.Range("D4").Select()

# --- Classes sheet ---
$classes = $wb.Worksheets.Item("Classes")

# Add new class row
$classes.Range("A6").Value = "STAT 102"
$classes.Range("B6").Value = 3

$classes.Range("C6").Select()

# Re-select the Teacher sheet so it remains the active tab
$teacher.Activate()
$teacher.Range("D4").Select()
